# Apply "BID API Completed - biding & Listing" edit:
# - The "Bidding API" row (row 14) on the "To Do" sheet is marked Done.
# - Column B, row 14 changes from "++" to "Done" (keeping the quote-prefix
#   formatting the cell already had, via a leading apostrophe).
# - The selected range moves from B14 to A14:B14 (active cell A14).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("To Do")

# Update status cell B14 from "++" to "Done"
$ws.Range("B14").Value = "'Done"

# Update the sheet's current selection to A14:B14 with A14 active
$ws.Range("A14:B14").Select()
